$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 12.600043675037883
$ws.Range("C2").Value = 11.862081862994724
$ws.Range("D2").Value = 13.299095352106743
$ws.Range("E2").Value = 12.944287081014441

$ws.Range("B3").Value = 12.855681704249111
$ws.Range("C3").Value = 10.818102188479116
$ws.Range("D3").Value = 14.134549322159225
$ws.Range("E3").Value = 10.789385201668139

$ws.Range("B1:E3").Select()
